# This script applies the odds/value updates described by the commit diff
# ("Atualizando o arquivo XLSX") to the single worksheet in the workbook.
# Only numeric <v> values change; no rows/columns are added or removed and
# no formatting/styles are touched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.02
$ws.Range("H2").Value = 1.02
$ws.Range("J2").Value = 1.02
$ws.Range("L2").Value = 1.01
$ws.Range("M2").Value = 1.01
$ws.Range("N2").Value = 1.25
$ws.Range("O2").Value = 1.01
$ws.Range("P2").Value = 1.25
$ws.Range("R2").Value = 1.18
$ws.Range("S2").Value = 1.4
$ws.Range("T2").Value = 1.01
$ws.Range("U2").Value = 1.01
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1000
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3
$ws.Range("F3").Value = 1.8
$ws.Range("G3").Value = 1.81
$ws.Range("H3").Value = 5.1
$ws.Range("J3").Value = 4
$ws.Range("Q3").Value = 1.6

# Row 4
$ws.Range("P4").Value = 1.5
$ws.Range("Q4").Value = 2.38

# Row 5
$ws.Range("F5").Value = 6.6
$ws.Range("G5").Value = 9.800000000000001
$ws.Range("H5").Value = 1.5
$ws.Range("I5").Value = 1.62
$ws.Range("K5").Value = 4.6
$ws.Range("P5").Value = 1.8
$ws.Range("Q5").Value = 2.02

# Row 9
$ws.Range("G9").Value = 3.8
$ws.Range("I9").Value = 2.22
$ws.Range("J9").Value = 3.5
$ws.Range("K9").Value = 3.75

# Row 10
$ws.Range("F10").Value = 3.05
$ws.Range("G10").Value = 5.2
$ws.Range("H10").Value = 1.81
$ws.Range("I10").Value = 2.18
$ws.Range("J10").Value = 3.4
$ws.Range("K10").Value = 5.4
$ws.Range("P10").Value = 2.22
$ws.Range("Q10").Value = 1.59

# Row 11
$ws.Range("F11").Value = 1.24
$ws.Range("G11").Value = 1.36
$ws.Range("H11").Value = 13
$ws.Range("I11").Value = 18.5
$ws.Range("K11").Value = 7.6

# Row 14
$ws.Range("I14").Value = 1.38

# Row 16
$ws.Range("F16").Value = 2.62
$ws.Range("G16").Value = 3.2

# Row 17
$ws.Range("F17").Value = 1.67
$ws.Range("G17").Value = 1.86
$ws.Range("H17").Value = 5.8
$ws.Range("I17").Value = 8
$ws.Range("J17").Value = 3.45
$ws.Range("K17").Value = 4.1
$ws.Range("P17").Value = 1.67

# Row 18
$ws.Range("H18").Value = 3.1
$ws.Range("U18").Value = 2.12
$ws.Range("AD18").Value = 13

# Row 19
$ws.Range("F19").Value = 1.84
$ws.Range("G19").Value = 1.85
$ws.Range("I19").Value = 5.5
$ws.Range("K19").Value = 3.7
$ws.Range("R19").Value = 1.24
$ws.Range("AL19").Value = 55

# Row 20
$ws.Range("H20").Value = 10.5
$ws.Range("I20").Value = 11.5
$ws.Range("P20").Value = 2.16
$ws.Range("Q20").Value = 1.75
